# feat: add 2022-Q4 data
#
# - Inserts a new "2022-Q4" sheet right after "总计" (i.e. before the
#   existing "2022-Q3" sheet), re-using the "2022-Q3" sheet's layout /
#   formatting as a starting point (same header row + index-column style).
# - Fills the new sheet with the 2022-Q4 fund-holding detail rows.
# - Adds a 2022-Q4 summary row to the "总计" sheet, shifting the existing
#   2022-Q3 / 2022-Q2 summary rows down by one row.

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item(1)
$q3Sheet = $wb.Worksheets.Item(2)

# --- 1. Create the new "2022-Q4" sheet right after "总计" -----------------
#        by duplicating "2022-Q3" (keeps sheetPr/margins/header style) ----
$q3Sheet.Copy($null, $totalSheet)
$q4Sheet = $wb.Worksheets.Item(2)
$q4Sheet.Name = "2022-Q4"

# The duplicated sheet only has 3 data rows (rows 2-4); the Q4 data needs
# 6 rows, so extend the formatted index-column / row style down to row 7.
$q4Sheet.Range("A4:H4").Copy()
$q4Sheet.Range("A5:H7").PasteSpecial(-4122)

# --- 2. Populate the "2022-Q4" sheet data ---------------------------------
$q4Data = @(
    @(0, "011603", "兴业高端制造混合A",          "0.61", "89.27", "3.28", "0.0200", 10),
    @(1, "011604", "兴业高端制造混合C",          "0.51", "89.27", "3.28", "0.0167", 10),
    @(2, "016648", "兴业数字经济优选股票C",        "0.44", "48.54", "3.67", "0.0161", 3),
    @(3, "001535", "景顺长城改革机遇灵活配置混合A", "0.28", "43.79", "1.47", "0.0041", 10),
    @(4, "007945", "景顺长城改革机遇灵活配置混合C", "0.24", "43.79", "1.47", "0.0035", 10),
    @(5, "016647", "兴业数字经济优选股票A",        "0.02", "48.54", "3.67", "0.0007", 3)
)

# B-G columns hold text (even the number-looking ones, e.g. "0.0200" must
# keep its trailing zero), so force the range to text before writing.
$q4Sheet.Range("B2:G7").NumberFormat = "@"

$r = 2
foreach ($row in $q4Data) {
    $q4Sheet.Range("A$r").Value = $row[0]
    $q4Sheet.Range("B$r").Value = $row[1]
    $q4Sheet.Range("C$r").Value = $row[2]
    $q4Sheet.Range("D$r").Value = $row[3]
    $q4Sheet.Range("E$r").Value = $row[4]
    $q4Sheet.Range("F$r").Value = $row[5]
    $q4Sheet.Range("G$r").Value = $row[6]
    $q4Sheet.Range("H$r").Value = $row[7]
    $r = $r + 1
}

# --- 3. Update the "总计" sheet: push the Q3 / Q2 rows down one row and ---
#        insert the new Q4 summary row at row 2 -----------------------------

# Extend the formatted index-column style down to the new row 4 first.
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A4").PasteSpecial(-4122)

$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2022-Q2"
$totalSheet.Range("C4").Value = 2
$totalSheet.Range("D4").Value = 0.8100000000000001

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2022-Q3"
$totalSheet.Range("C3").Value = 3
$totalSheet.Range("D3").Value = 0.07000000000000001

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q4"
$totalSheet.Range("C2").Value = 6
$totalSheet.Range("D2").Value = 0.06
